$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 197.97
$ws.Range("I15").Value = 197.97
$ws.Range("K15").Value = 593.91
$ws.Range("M15").Value = -424.91
$ws.Range("H98").Value = 1960.5349
$ws.Range("I98").Value = 1961.641
$ws.Range("J98").Value = 1949.75
$ws.Range("K98").Value = 1961.641
$ws.Range("L98").Value = 1949.75
$ws.Range("M98").Value = -463.6410000000001
$ws.Range("N98").Value = -4945.75
$ws.Range("H122").Value = 1960.5349
$ws.Range("I122").Value = 1961.641
$ws.Range("J122").Value = 1949.75
$ws.Range("K122").Value = 5884.923000000001
$ws.Range("L122").Value = 5849.25
$ws.Range("M122").Value = -3434.923000000001
$ws.Range("N122").Value = -10749.25
$ws.Range("H127").Value = 100001260
$ws.Range("I127").Value = 1225
$ws.Range("J127").Value = 166667950
$ws.Range("K127").Value = 3675
$ws.Range("L127").Value = 500003850
$ws.Range("M127").Value = 1285
$ws.Range("N127").Value = -500013770
$ws.Range("H138").Value = 3827.3691
$ws.Range("J138").Value = 4765.8213
$ws.Range("L138").Value = 14297.4639
$ws.Range("N138").Value = -24577.4639
$ws.Range("H139").Value = 68619
$ws.Range("J139").Value = 68619
$ws.Range("L139").Value = 68619
$ws.Range("N139").Value = -78899
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1696.4054
$ws.Range("I2").Value = 1832.8695
$ws.Range("J2").Value = 1472.2142
$ws.Range("K2").Value = 1832.8695
$ws.Range("L2").Value = 1472.2142
$ws.Range("M2").Value = -1719.8695
$ws.Range("N2").Value = -1698.2142
$ws.Range("H74").Value = 97534.42
$ws.Range("I74").Value = 108862.336
$ws.Range("K74").Value = 108862.336
$ws.Range("M74").Value = -107988.336
$ws.Range("H77").Value = 97534.42
$ws.Range("I77").Value = 108862.336
$ws.Range("K77").Value = 544311.6799999999
$ws.Range("M77").Value = -539943.6799999999
$ws.Range("H116").Value = 1696.4054
$ws.Range("I116").Value = 1832.8695
$ws.Range("J116").Value = 1472.2142
$ws.Range("K116").Value = 1832.8695
$ws.Range("L116").Value = 1472.2142
$ws.Range("M116").Value = 461.1305
$ws.Range("N116").Value = -6060.2142
$ws.Range("H122").Value = 7814662.5
$ws.Range("I122").Value = 2133.3333
$ws.Range("K122").Value = 6399.999899999999
$ws.Range("M122").Value = -3949.999899999999
$ws.Range("H133").Value = 51163.65
$ws.Range("J133").Value = 51163.65
$ws.Range("L133").Value = 51163.65
$ws.Range("N133").Value = -56223.65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1696.4054
$ws.Range("I3").Value = 1832.8695
$ws.Range("J3").Value = 1472.2142
$ws.Range("K3").Value = 1832.8695
$ws.Range("L3").Value = 1472.2142
$ws.Range("M3").Value = -1718.8695
$ws.Range("N3").Value = -1700.2142
$ws.Range("H20").Value = 939.15
$ws.Range("I20").Value = 884.93335
$ws.Range("J20").Value = 1101.8
$ws.Range("K20").Value = 884.93335
$ws.Range("L20").Value = 1101.8
$ws.Range("M20").Value = -637.93335
$ws.Range("N20").Value = -1595.8
$ws.Range("H33").Value = 22999
$ws.Range("I33").Value = 22999
$ws.Range("K33").Value = 22999
$ws.Range("M33").Value = -22663
$ws.Range("H132").Value = 59356.957
$ws.Range("J132").Value = 59356.957
$ws.Range("L132").Value = 59356.957
$ws.Range("N132").Value = -69476.95699999999
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
$ws.Range("H134").Value = 3402.6072
$ws.Range("I134").Value = 3375.3462
$ws.Range("K134").Value = 10126.0386
$ws.Range("M134").Value = -7591.0386
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1563.6666
$ws.Range("I16").Value = 1595.5
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1595.5
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1308.5
$ws.Range("N16").Value = -2074
$ws.Range("H99").Value = 3800
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 1563.6666
$ws.Range("I113").Value = 1595.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1595.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 574.5
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 6343.7617
$ws.Range("I122").Value = 5018.727
$ws.Range("K122").Value = 15056.181
$ws.Range("M122").Value = -12606.181
$ws.Range("H126").Value = 3800
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15629163
$ws.Range("I5").Value = 705.2727
$ws.Range("J5").Value = 50011772
$ws.Range("K5").Value = 2115.8181
$ws.Range("L5").Value = 150035316
$ws.Range("M5").Value = -2003.8181
$ws.Range("N5").Value = -150035540
$ws.Range("H12").Value = 32258282
$ws.Range("J12").Value = 218.1579
$ws.Range("L12").Value = 654.4737
$ws.Range("N12").Value = -1000.4737
$ws.Range("H40").Value = 97.59999999999999
$ws.Range("I40").Value = 97.59999999999999
$ws.Range("K40").Value = 390.4
$ws.Range("M40").Value = -321.4
$ws.Range("H110").Value = 4997
$ws.Range("I110").Value = 1900
$ws.Range("J110").Value = 6029.3335
$ws.Range("K110").Value = 5700
$ws.Range("L110").Value = 18088.0005
$ws.Range("M110").Value = -1610
$ws.Range("N110").Value = -26268.0005
$ws.Range("H131").Value = 42289.043
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 42289.043
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 126867.129
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -136947.129
$ws.Range("H135").Value = 15629163
$ws.Range("I135").Value = 705.2727
$ws.Range("J135").Value = 50011772
$ws.Range("K135").Value = 6347.454299999999
$ws.Range("L135").Value = 450105948
$ws.Range("M135").Value = -3812.454299999999
$ws.Range("N135").Value = -450111018
$ws.Range("H137").Value = 71433640
$ws.Range("I137").Value = 125003624
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 375010872
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -375005772
$ws.Range("N137").Value = -31200
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 27971.428
$ws.Range("J123").Value = 27971.428
$ws.Range("L123").Value = 27971.428
$ws.Range("N123").Value = -32871.428
$ws.Range("H126").Value = 2595.0356
$ws.Range("I126").Value = 1848
$ws.Range("J126").Value = 3457
$ws.Range("K126").Value = 5544
$ws.Range("L126").Value = 10371
$ws.Range("M126").Value = -3074
$ws.Range("N126").Value = -15311
$ws.Range("H132").Value = 25718.582
$ws.Range("I132").Value = 60329.47
$ws.Range("K132").Value = 180988.41
$ws.Range("M132").Value = -178458.41
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 20001
$ws.Range("J13").Value = 20001
$ws.Range("L13").Value = 20001
$ws.Range("N13").Value = -20281
$ws.Range("H40").Value = 3691.6216
$ws.Range("I40").Value = 3565.3845
$ws.Range("J40").Value = 3990
$ws.Range("K40").Value = 3565.3845
$ws.Range("L40").Value = 3990
$ws.Range("M40").Value = -3429.3845
$ws.Range("N40").Value = -4262
$ws.Range("H58").Value = 34500
$ws.Range("J58").Value = 34500
$ws.Range("L58").Value = 34500
$ws.Range("N58").Value = -35020
$ws.Range("H93").Value = 1936
$ws.Range("I93").Value = 1451.25
$ws.Range("K93").Value = 1451.25
$ws.Range("M93").Value = -203.25
$ws.Range("H122").Value = 5838.8945
$ws.Range("I122").Value = 4652.7856
$ws.Range("J122").Value = 9160
$ws.Range("K122").Value = 13958.3568
$ws.Range("L122").Value = 27480
$ws.Range("M122").Value = -11508.3568
$ws.Range("N122").Value = -32380
$ws.Range("H132").Value = 2718.6667
$ws.Range("I132").Value = 2117.32
$ws.Range("J132").Value = 4085.3635
$ws.Range("K132").Value = 6351.960000000001
$ws.Range("L132").Value = 12256.0905
$ws.Range("M132").Value = -3821.960000000001
$ws.Range("N132").Value = -17316.0905
$ws.Range("H136").Value = 4161.078
$ws.Range("I136").Value = 2560.3572
$ws.Range("J136").Value = 7217
$ws.Range("K136").Value = 7681.071599999999
$ws.Range("L136").Value = 21651
$ws.Range("M136").Value = -5131.071599999999
$ws.Range("N136").Value = -26751
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 800
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H139").Value = 69480.71000000001
$ws.Range("J139").Value = 69480.71000000001
$ws.Range("L139").Value = 69480.71000000001
$ws.Range("N139").Value = -79760.71000000001
